$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.794.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "'1.628.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'215.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'0.5069"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.2580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.06443"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "'0.07789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'4.258"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "'1.852.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'1.623.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "'0.5569"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "'63.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").Value = "'0.0₅7547"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "'25.796.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'193.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "'4.298"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("D22").Value = "'9.813"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "'6.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'1.806"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.55%  "
$ws.Range("D26").Value = "'140.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'6.724"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'0.04868"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "'3.275"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'3.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "'1.555"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "'2.371"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'0.8936"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "'2.570"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'1.132.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("D39").Value = "'0.5460"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "'0.01554"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'5.565"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").Value = "'0.7953"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "'97.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "'1.780.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'0.0₈113"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("D47").Value = "'0.4431"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "'55.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "'0.05052"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "'1.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
